# Stats.xlsx update: add "2023-01" column (E) of results to the existing
# benchmark tables on the "Positive" sheet, and extend the corresponding
# bar charts with the new data series.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Positive")

# --- Table 1: rows 1-16 ("Sorting 10 million ints ... 0 to 10 millions") ---
# No 2023-01 figures were actually captured for this table; the author
# only left an (empty, quote-prefixed) placeholder header cell.
$ws.Range("E2").Value = "'"

# --- Table 2: rows 25-40 ("Sorting 10 million ints ... 0 to 100000") ---
$ws.Range("E26").Value = "'2023-01"
$ws.Range("E27").Value = 524
$ws.Range("E28").Value = 44
$ws.Range("E29").Value = 75
$ws.Range("E30").Value = 90
$ws.Range("E31").Value = 72
$ws.Range("E32").Value = 45
$ws.Range("E33").Value = 39
$ws.Range("E34").Value = 41

$ws.Range("E36").Formula = "=E26"
$ws.Range("E37").Formula = "=E31"
$ws.Range("E38").Formula = "=E32"
$ws.Range("E39").Formula = "=E33"
$ws.Range("E40").Formula = "=E34"

# --- Table 3: rows 50-65 ("Sorting 40 million ints ... 0 to 1000 Millions") ---
$ws.Range("E51").Value = "'2023-01"
$ws.Range("E52").Value = 3096
$ws.Range("E53").Value = 2486
$ws.Range("E54").Value = 433
$ws.Range("E55").Value = 439
$ws.Range("E56").Value = 353
$ws.Range("E57").Value = 583
$ws.Range("E58").Value = 379
$ws.Range("E59").Value = 252

$ws.Range("E61").Formula = "=E51"
$ws.Range("E62").Formula = "=E56"
$ws.Range("E63").Formula = "=E57"
$ws.Range("E64").Formula = "=E58"
$ws.Range("E65").Formula = "=E59"

# --- Table 4: rows 75-80 ("Sorting 10 million objects ... 0 to 10 Millions") ---
$ws.Range("E77").Value = "'2023-01"
$ws.Range("E78").Value = 5416
$ws.Range("E79").Value = 691
$ws.Range("E80").Value = 965

# --- Table 5: rows 86-91 ("Sorting 10 million objects ... 0 to 100000") ---
$ws.Range("E88").Value = "'2023-01"
$ws.Range("E89").Value = 4089
$ws.Range("E90").Value = 637
$ws.Range("E91").Value = 608

# Restore the selection the author left behind after entering the data.
$ws.Range("E81").Select()
